# Sun, Apr 05, 2020  7:06:33 PM
#
# 1) Slide 16's table (3rd shape on the slide) switches from the deck's
#    one custom table style to a built-in PowerPoint table style.
# 2) The presentation's theme color scheme (ppt/theme/theme2.xml, the
#    theme actually wired to the slide master / used by the deck) is
#    swapped from the "Integral" palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 ------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{349A6D5B-C2BA-4A32-BB42-47B488E3DABF}")

# --- 2. Theme color scheme: Integral -> Office ---------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# .RGB takes a PowerPoint-style color long (R + G*256 + B*65536), i.e. the
# hex value with its bytes reversed from "RRGGBB". MsoThemeColorSchemeIndex
# order is: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$themeColors.Item(1).RGB  = 0x000000    # dk1      -> 000000
$themeColors.Item(2).RGB  = 0xFFFFFF    # lt1      -> FFFFFF
$themeColors.Item(3).RGB  = 0x6A5444    # dk2      -> 44546A
$themeColors.Item(4).RGB  = 0xE6E6E7    # lt2      -> E7E6E6
$themeColors.Item(5).RGB  = 0xD59B5B    # accent1  -> 5B9BD5
$themeColors.Item(6).RGB  = 0x317DED    # accent2  -> ED7D31
$themeColors.Item(7).RGB  = 0xA5A5A5    # accent3  -> A5A5A5
$themeColors.Item(8).RGB  = 0x00C0FF    # accent4  -> FFC000
$themeColors.Item(9).RGB  = 0xC47244    # accent5  -> 4472C4
$themeColors.Item(10).RGB = 0x47AD70    # accent6  -> 70AD47
$themeColors.Item(11).RGB = 0xC16305    # hlink    -> 0563C1
$themeColors.Item(12).RGB = 0x724F95    # folHlink -> 954F72
